$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.86%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "34.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.18%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.197"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.40%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07821"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.68%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.379"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.45%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.032"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.09%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.960"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "6.29%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9330"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.36%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1015"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9.41%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1853"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "9.26%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08596"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.60%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03315"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.28%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09897"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.49%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001488"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.45%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005718"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.13%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.471"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.76%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.149"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.02%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3388"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.73%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1303"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.31%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.306"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.51%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2221"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.68%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04574"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.95%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001216"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.55%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004441"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.68%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001296"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.39%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003400"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.13%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "13.17%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04797"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "7.99%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007752"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.57%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1408"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.006938"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-27.27%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002203"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.84%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009465"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.15%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005908"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.39%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.36%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "22.93%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002094"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.36%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001994"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.36%"
